$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 196 ----
# Column A keeps the same date/time style as the rest of the column, so
# copy the format from the last existing data row before writing the value.
$ws.Range("A195").Copy() | Out-Null
$ws.Range("A196").PasteSpecial(-4122) | Out-Null
$ws.Range("A196").Value = 45477.2916666667

$ws.Range("B196").Value = 27000
$ws.Range("C196").Value = 3.77999997138977
$ws.Range("D196").Value = 3.77999997138977
$ws.Range("E196").Value = 3.77999997138977
$ws.Range("F196").Value = 3.77999997138977

# Column G holds the "close" value again, but stored as text (shared
# string), matching the source data. Force text storage for this
# numeric-looking string, then drop the format override again so the
# cell is left with the default (unstyled) look, same as its peers.
$ws.Range("G196").NumberFormat = "@"
$ws.Range("G196").Value = "3.77999997138977"
$ws.Range("G196").ClearFormats()

$ws.Range("H196").Value = "ELSA.MI"

# ---- Row 197 ----
$ws.Range("A195").Copy() | Out-Null
$ws.Range("A197").PasteSpecial(-4122) | Out-Null
$ws.Range("A197").Value = 45478.4490740741

$ws.Range("B197").Value = 500
$ws.Range("C197").Value = 3.8199999332428
$ws.Range("D197").Value = 3.8199999332428
$ws.Range("E197").Value = 3.8199999332428
$ws.Range("F197").Value = 3.8199999332428

$ws.Range("G197").NumberFormat = "@"
$ws.Range("G197").Value = "3.8199999332428"
$ws.Range("G197").ClearFormats()

$ws.Range("H197").Value = "ELSA.MI"
